$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value2 = 3823001401
$ws.Range("C3").Value2 = 90
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "70005873"
$ws.Range("D3").Style = "Normal"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value2 = "ZNPL"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value2 = 1600
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value2 = "20220319"
$ws.Range("J3").Style = "Normal"

# Row 4
$ws.Range("B4").Value2 = 3823001401
$ws.Range("C4").Value2 = 90
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "70013970"
$ws.Range("D4").Style = "Normal"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value2 = "ZCQE"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value2 = 3200
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value2 = "20220319"
$ws.Range("J4").Style = "Normal"

# Row 5
$ws.Range("B5").Value2 = 3820017820
$ws.Range("C5").Value2 = 30
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "70003132"
$ws.Range("D5").Style = "Normal"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value2 = "ZNPL"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value2 = 64
$ws.Range("J5").NumberFormat = "@"
$ws.Range("J5").Value2 = "20220226"
$ws.Range("J5").Style = "Normal"

# Row 6
$ws.Range("B6").Value2 = 3820017820
$ws.Range("C6").Value2 = 30
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "70026583"
$ws.Range("D6").Style = "Normal"
$ws.Range("F6").Value2 = ""
$ws.Range("G6").Value2 = 500
$ws.Range("J6").NumberFormat = "@"
$ws.Range("J6").Value2 = "20220226"
$ws.Range("J6").Style = "Normal"

# Row 7
$ws.Range("B7").Value2 = 3820017820
$ws.Range("C7").Value2 = 30
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "70011454"
$ws.Range("D7").Style = "Normal"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value2 = "ZCMQ"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value2 = 500
$ws.Range("J7").NumberFormat = "@"
$ws.Range("J7").Value2 = "20220226"
$ws.Range("J7").Style = "Normal"

# Row 8
$ws.Range("B8").Value2 = 3820017820
$ws.Range("C8").Value2 = 30
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "72000192"
$ws.Range("D8").Style = "Normal"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value2 = "ZCMQ"
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value2 = 1000
$ws.Range("J8").NumberFormat = "@"
$ws.Range("J8").Value2 = "20220226"
$ws.Range("J8").Style = "Normal"

# Row 9
$ws.Range("B9").Value2 = 3820017820
$ws.Range("C9").Value2 = 30
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "70012107"
$ws.Range("D9").Style = "Normal"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value2 = "ZCMQ"
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Value2 = 1890
$ws.Range("J9").NumberFormat = "@"
$ws.Range("J9").Value2 = "20220226"
$ws.Range("J9").Style = "Normal"

# Row 10
$ws.Range("B10").Value2 = 3820017820
$ws.Range("C10").Value2 = 30
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "72000211"
$ws.Range("D10").Style = "Normal"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value2 = "ZCMQ"
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").Value2 = 1000
$ws.Range("J10").NumberFormat = "@"
$ws.Range("J10").Value2 = "20220226"
$ws.Range("J10").Style = "Normal"

# Row 11
$ws.Range("B11").Value2 = 3820017820
$ws.Range("C11").Value2 = 30
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "72002008"
$ws.Range("D11").Style = "Normal"
$ws.Range("G11").Value2 = 1000
$ws.Range("J11").NumberFormat = "@"
$ws.Range("J11").Value2 = "20220226"
$ws.Range("J11").Style = "Normal"

# Row 12
$ws.Range("B12").Value2 = 3820017820
$ws.Range("C12").Value2 = 30
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "72000768"
$ws.Range("D12").Style = "Normal"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value2 = "ZCEQ"
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").Value2 = 1000
$ws.Range("J12").NumberFormat = "@"
$ws.Range("J12").Value2 = "20220226"
$ws.Range("J12").Style = "Normal"

# Row 13
$ws.Range("B13").Value2 = 3820017820
$ws.Range("C13").Value2 = 20
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "70003132"
$ws.Range("D13").Style = "Normal"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value2 = "ZNPL"
$ws.Range("F13").Style = "Normal"
$ws.Range("G13").Value2 = 39
$ws.Range("J13").NumberFormat = "@"
$ws.Range("J13").Value2 = "20220226"
$ws.Range("J13").Style = "Normal"

# Row 14
$ws.Range("B14").Value2 = 3820017820
$ws.Range("C14").Value2 = 20
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "70026583"
$ws.Range("D14").Style = "Normal"
$ws.Range("F14").Value2 = ""
$ws.Range("G14").Value2 = 201
$ws.Range("J14").NumberFormat = "@"
$ws.Range("J14").Value2 = "20220226"
$ws.Range("J14").Style = "Normal"

# Row 15
$ws.Range("B15").Value2 = 3820017820
$ws.Range("C15").Value2 = 20
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "70011454"
$ws.Range("D15").Style = "Normal"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value2 = "ZCMQ"
$ws.Range("F15").Style = "Normal"
$ws.Range("G15").Value2 = 201
$ws.Range("J15").NumberFormat = "@"
$ws.Range("J15").Value2 = "20220226"
$ws.Range("J15").Style = "Normal"

# Row 16
$ws.Range("B16").Value2 = 3820017820
$ws.Range("C16").Value2 = 20
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "72000128"
$ws.Range("D16").Style = "Normal"
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value2 = "ZCMQ"
$ws.Range("F16").Style = "Normal"
$ws.Range("G16").Value2 = 1204
$ws.Range("J16").NumberFormat = "@"
$ws.Range("J16").Value2 = "20220226"
$ws.Range("J16").Style = "Normal"

# Row 17
$ws.Range("B17").Value2 = 3820017820
$ws.Range("C17").Value2 = 20
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "70012107"
$ws.Range("D17").Style = "Normal"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value2 = "ZCMQ"
$ws.Range("F17").Style = "Normal"
$ws.Range("G17").Value2 = 1166
$ws.Range("J17").NumberFormat = "@"
$ws.Range("J17").Value2 = "20220226"
$ws.Range("J17").Style = "Normal"

# Row 18
$ws.Range("B18").Value2 = 3820017820
$ws.Range("C18").Value2 = 20
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "72000191"
$ws.Range("D18").Style = "Normal"
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value2 = "ZCMQ"
$ws.Range("F18").Style = "Normal"
$ws.Range("G18").Value2 = 1204
$ws.Range("J18").NumberFormat = "@"
$ws.Range("J18").Value2 = "20220226"
$ws.Range("J18").Style = "Normal"

# Row 19
$ws.Range("B19").Value2 = 3820017820
$ws.Range("C19").Value2 = 20
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "72002007"
$ws.Range("D19").Style = "Normal"
$ws.Range("F19").Value2 = ""
$ws.Range("G19").Value2 = 1204
$ws.Range("J19").NumberFormat = "@"
$ws.Range("J19").Value2 = "20220226"
$ws.Range("J19").Style = "Normal"

# Row 20
$ws.Range("B20").Value2 = 3820017820
$ws.Range("C20").Value2 = 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "72000764"
$ws.Range("D20").Style = "Normal"
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value2 = "ZCEQ"
$ws.Range("F20").Style = "Normal"
$ws.Range("G20").Value2 = 1204
$ws.Range("J20").NumberFormat = "@"
$ws.Range("J20").Value2 = "20220226"
$ws.Range("J20").Style = "Normal"

# Row 21
$ws.Range("B21").Value2 = 3820017820
$ws.Range("C21").Value2 = 10
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "70003132"
$ws.Range("D21").Style = "Normal"
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value2 = "ZNPL"
$ws.Range("F21").Style = "Normal"
$ws.Range("G21").Value2 = 81
$ws.Range("J21").NumberFormat = "@"
$ws.Range("J21").Value2 = "20220122"
$ws.Range("J21").Style = "Normal"

# Row 22
$ws.Range("B22").Value2 = 3820017820
$ws.Range("C22").Value2 = 10
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "70026583"
$ws.Range("D22").Style = "Normal"
$ws.Range("F22").Value2 = ""
$ws.Range("G22").Value2 = 268
$ws.Range("J22").NumberFormat = "@"
$ws.Range("J22").Value2 = "20220122"
$ws.Range("J22").Style = "Normal"

# Row 23
$ws.Range("B23").Value2 = 3820017820
$ws.Range("C23").Value2 = 10
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "70011454"
$ws.Range("D23").Style = "Normal"
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value2 = "ZCMQ"
$ws.Range("F23").Style = "Normal"
$ws.Range("G23").Value2 = 268
$ws.Range("J23").NumberFormat = "@"
$ws.Range("J23").Value2 = "20220122"
$ws.Range("J23").Style = "Normal"

# Row 24 (new)
$ws.Range("A24").Value2 = 22
$ws.Range("B24").Value2 = 3820017820
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value2 = 10
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "72001086"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = 3801
$ws.Range("E24").Style = "Normal"
$ws.Range("F24").NumberFormat = "@"
$ws.Range("F24").Value2 = "ZCMQ"
$ws.Range("F24").Style = "Normal"
$ws.Range("G24").Value2 = 2412
$ws.Range("G24").Style = "Normal"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value2 = "0001"
$ws.Range("H24").Style = "Normal"
$ws.Range("J24").NumberFormat = "@"
$ws.Range("J24").Value2 = "20220212"
$ws.Range("J24").Style = "Normal"
$ws.Range("M24").Value2 = "X"
$ws.Range("N24").Value2 = "X"

# Row 25 (new)
$ws.Range("A25").Value2 = 23
$ws.Range("B25").Value2 = 3820017820
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value2 = 10
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "70012107"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = 3801
$ws.Range("E25").Style = "Normal"
$ws.Range("F25").NumberFormat = "@"
$ws.Range("F25").Value2 = "ZCMQ"
$ws.Range("F25").Style = "Normal"
$ws.Range("G25").Value2 = 2412
$ws.Range("G25").Style = "Normal"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value2 = "0001"
$ws.Range("H25").Style = "Normal"
$ws.Range("J25").NumberFormat = "@"
$ws.Range("J25").Value2 = "20220129"
$ws.Range("J25").Style = "Normal"
$ws.Range("M25").Value2 = "X"
$ws.Range("N25").Value2 = "X"

# Row 26 (new)
$ws.Range("A26").Value2 = 24
$ws.Range("B26").Value2 = 3820017820
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value2 = 10
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "72000212"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = 3801
$ws.Range("E26").Style = "Normal"
$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value2 = "ZCMQ"
$ws.Range("F26").Style = "Normal"
$ws.Range("G26").Value2 = 2412
$ws.Range("G26").Style = "Normal"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value2 = "0001"
$ws.Range("H26").Style = "Normal"
$ws.Range("J26").NumberFormat = "@"
$ws.Range("J26").Value2 = "20220129"
$ws.Range("J26").Style = "Normal"
$ws.Range("M26").Value2 = "X"
$ws.Range("N26").Value2 = "X"

# Row 27 (new)
$ws.Range("A27").Value2 = 25
$ws.Range("B27").Value2 = 3820017820
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value2 = 10
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "72002006"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = 3801
$ws.Range("E27").Style = "Normal"
$ws.Range("F27").Value2 = ""
$ws.Range("F27").Style = "Normal"
$ws.Range("G27").Value2 = 2412
$ws.Range("G27").Style = "Normal"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value2 = "0001"
$ws.Range("H27").Style = "Normal"
$ws.Range("J27").NumberFormat = "@"
$ws.Range("J27").Value2 = "20220226"
$ws.Range("J27").Style = "Normal"
$ws.Range("M27").Value2 = "X"
$ws.Range("N27").Value2 = "X"

# Row 28 (new)
$ws.Range("A28").Value2 = 26
$ws.Range("B28").Value2 = 3820017820
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value2 = 10
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "72000763"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = 3801
$ws.Range("E28").Style = "Normal"
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value2 = "ZCEQ"
$ws.Range("F28").Style = "Normal"
$ws.Range("G28").Value2 = 2412
$ws.Range("G28").Style = "Normal"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value2 = "0001"
$ws.Range("H28").Style = "Normal"
$ws.Range("J28").NumberFormat = "@"
$ws.Range("J28").Value2 = "20220226"
$ws.Range("J28").Style = "Normal"
$ws.Range("M28").Value2 = "X"
$ws.Range("N28").Value2 = "X"
